$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update logged hours for Tues (row3) and Wed (row4) in the "Wed"/"Tues" test-motion column (C)
$ws.Range("C3").Value = 4.5
$ws.Range("C4").Value = 2.5

# Update the active selection to C5
$ws.Range("C5").Select()
